$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The forecast for 3/25 (row 20) is no longer a projection once the actual
# confirmed-case count for that day is known, so replace the forecast
# formula in I20 with the real observed value. First copy the formatting
# used by the other "actual" (non-forecast) cells in column I (e.g. I19)
# onto I20, then overwrite the cell's content with the hard value.
$ws.Range("I19").Copy()
$ws.Range("I20").PasteSpecial(-4122)
$ws.Range("I20").Value = 68211

# Update the active selection to I21, matching the saved view state.
$ws.Range("I21").Select()
